$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

# Target DAS values (Squaring, Flowering, Cutout, Openbolls, Maturity, HarvestRipe)
# for Phenology columns R:W, applied uniformly to rows 2-22.
$cols = @("R", "S", "T", "U", "V", "W")
$values = @(32, 48, 83, 111, 141, 160)

# Reference cells that already carry the correct number format/style for
# each column (sampled from the pre-existing seed values in the sheet) so
# that filling in the remaining rows reuses the same style rather than
# minting new ones.
$styleRefs = @{
    "R" = "R6"
    "S" = "S9"
    "T" = "T14"
    "U" = "U18"
    "V" = "V20"
    "W" = "R6"
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $val = $values[$i]
    $refCell = $ws.Range($styleRefs[$col])
    $refCell.Copy()
    for ($row = 2; $row -le 22; $row++) {
        $target = $ws.Range($col + $row)
        $target.PasteSpecial(-4122)
        $target.Value = $val
    }
}
